$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving plain numeric-looking text must be forced to Text format
# so Excel does not auto-convert them to numbers, matching the original inlineStr type.
$textCells = @("D5", "D9", "D10", "D11", "D17", "D20", "D21", "D22", "D25", "D28", "D30", "D45", "D46", "D48", "D49")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.391.04'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '1.626.08'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '212.66'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').Value = '0.0618'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = '18.93'
$ws.Range('D11').Value = '0.0815'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').Value = '1.851.70'
$ws.Range('E12').Value = '  +1.50%  '
$ws.Range('D13').Value = '1.625.13'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').Value = '26.394.78'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '62.67'
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '202.98'
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('D21').Value = '4.28'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = '9.33'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('E24').Value = '  -3.41%  '
$ws.Range('D25').Value = '144.58'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -2.99%  '
$ws.Range('D28').Value = '15.22'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('D30').Value = '0.0517'
$ws.Range('E30').Value = '  +5.36%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  +1.66%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('D36').Value = '1.159.26'
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('E37').Value = '  +0.38%  '
$ws.Range('E38').Value = '  +2.68%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').Value = '1.763.22'
$ws.Range('D45').Value = '92.01'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').Value = '1.53'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('D47').Value = '0.0₆0104'
$ws.Range('E47').Value = '  +8.90%  '
$ws.Range('D48').Value = '54.05'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').Value = '0.0508'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('E51').Value = '  -0.17%  '

# Restore default (Normal) style so no stray number-format style attribute remains on the cells
foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}
